$d = $word.ActiveDocument

# --- Footers: both the default and first-page footer carry a Pearson
# Edexcel logo picture that was named "image1.png"; rename it to
# "image2.png" (the docPr/cNvPr shape name, not the underlying media file).
for ($k = 1; $k -le $d.Sections.Count; $k++) {
    $sec = $d.Sections($k)
    for ($i = 1; $i -le 3; $i++) {
        $ftr = $sec.Footers($i)
        if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -gt 0) {
            for ($j = 1; $j -le $ftr.Range.InlineShapes.Count; $j++) {
                $inlineShape = $ftr.Range.InlineShapes($j)
                if ($inlineShape.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shape = $inlineShape.ConvertToShape()
                    $shape.Name = "image2.png"
                    $shape.ConvertToInlineShape() | Out-Null
                }
            }
        }
    }

    # --- Headers: the first-page header carries the BTEC logo picture that
    # was named "image2.jpg"; rename it to "image1.jpg".
    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers($i)
        if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
            for ($j = 1; $j -le $hdr.Range.InlineShapes.Count; $j++) {
                $inlineShape = $hdr.Range.InlineShapes($j)
                if ($inlineShape.AlternativeText -eq "BTec_Logo-Orange") {
                    $shape = $inlineShape.ConvertToShape()
                    $shape.Name = "image1.jpg"
                    $shape.ConvertToInlineShape() | Out-Null
                }
            }
        }
    }
}

Write-Output "done"
